$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 36.81423266666667
$ws.Range("H2").Value = 110.442698
$ws.Range("I2").Value = 0.13776174071044
$ws.Range("J2").Value = 0.13776174071044
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 153.5290173333333
$ws.Range("N2").Value = 460.587052
$ws.Range("O2").Value = 0.3172206968818489
$ws.Range("P2").Value = 0.317220696881849
$ws.Range("Q2").Value = 5652.052965194032
$ws.Range("R2").Value = 50868.4766867463
$ws.Range("S2").Value = 0.04370087539182236
$ws.Range("T2").Value = 0.04370087539182237
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 36.81423266666667
$ws.Range("H3").Value = 110.442698
$ws.Range("I3").Value = 0.13776174071044
$ws.Range("J3").Value = 0.13776174071044
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3487728915577651
$ws.Range("P3").Value = 0.3487728915577651
$ws.Range("Q3").Value = 6214.231528034822
$ws.Range("R3").Value = 55928.0837523134
$ws.Range("S3").Value = 0.04804756065361125
$ws.Range("T3").Value = 0.04804756065361125
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 36.81423266666667
$ws.Range("H4").Value = 110.442698
$ws.Range("I4").Value = 0.13776174071044
$ws.Range("J4").Value = 0.13776174071044
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 68.09032333333333
$ws.Range("N4").Value = 204.27097
$ws.Range("O4").Value = 0.1406878008722904
$ws.Range("P4").Value = 0.1406878008722904
$ws.Range("Q4").Value = 2506.693005541896
$ws.Range("R4").Value = 22560.23704987706
$ws.Range("S4").Value = 0.01938139634489049
$ws.Range("T4").Value = 0.01938139634489049
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 36.81423266666667
$ws.Range("H5").Value = 110.442698
$ws.Range("I5").Value = 0.13776174071044
$ws.Range("J5").Value = 0.13776174071044
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 93.562673
$ws.Range("N5").Value = 280.688019
$ws.Range("O5").Value = 0.1933186106880956
$ws.Range("P5").Value = 0.1933186106880956
$ws.Range("Q5").Value = 3444.438012737251
$ws.Range("R5").Value = 30999.94211463527
$ws.Range("S5").Value = 0.02663190832011593
$ws.Range("T5").Value = 0.02663190832011593
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 214.101181
$ws.Range("H6").Value = 642.303543
$ws.Range("I6").Value = 0.8011833806175486
$ws.Range("J6").Value = 0.8011833806175486
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 153.5290173333333
$ws.Range("N6").Value = 460.587052
$ws.Range("O6").Value = 0.3172206968818489
$ws.Range("P6").Value = 0.317220696881849
$ws.Range("Q6").Value = 32870.74392883613
$ws.Range("R6").Value = 295836.6953595252
$ws.Range("S6").Value = 0.2541519503296543
$ws.Range("T6").Value = 0.2541519503296544
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 214.101181
$ws.Range("H7").Value = 642.303543
$ws.Range("I7").Value = 0.8011833806175486
$ws.Range("J7").Value = 0.8011833806175486
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3487728915577651
$ws.Range("P7").Value = 0.3487728915577651
$ws.Range("Q7").Value = 36140.21569338219
$ws.Range("R7").Value = 325261.9412404397
$ws.Range("S7").Value = 0.2794310443260079
$ws.Range("T7").Value = 0.2794310443260079
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 214.101181
$ws.Range("H8").Value = 642.303543
$ws.Range("I8").Value = 0.8011833806175486
$ws.Range("J8").Value = 0.8011833806175486
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 68.09032333333333
$ws.Range("N8").Value = 204.27097
$ws.Range("O8").Value = 0.1406878008722904
$ws.Range("P8").Value = 0.1406878008722904
$ws.Range("Q8").Value = 14578.21864033852
$ws.Range("R8").Value = 131203.9677630467
$ws.Range("S8").Value = 0.1127167279145101
$ws.Range("T8").Value = 0.1127167279145101
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 214.101181
$ws.Range("H9").Value = 642.303543
$ws.Range("I9").Value = 0.8011833806175486
$ws.Range("J9").Value = 0.8011833806175486
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 93.562673
$ws.Range("N9").Value = 280.688019
$ws.Range("O9").Value = 0.1933186106880956
$ws.Range("P9").Value = 0.1933186106880956
$ws.Range("Q9").Value = 20031.87878681681
$ws.Range("R9").Value = 180286.9090813513
$ws.Range("S9").Value = 0.1548836580473762
$ws.Range("T9").Value = 0.1548836580473762
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2781493333333334
$ws.Range("H10").Value = 0.8344480000000001
$ws.Range("I10").Value = 0.001040856580779521
$ws.Range("J10").Value = 0.001040856580779521
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 153.5290173333333
$ws.Range("N10").Value = 460.587052
$ws.Range("O10").Value = 0.3172206968818489
$ws.Range("P10").Value = 0.317220696881849
$ws.Range("Q10").Value = 42.70399381858844
$ws.Range("R10").Value = 384.335944367296
$ws.Range("S10").Value = 0.0003301812499089382
$ws.Range("T10").Value = 0.0003301812499089382
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2781493333333334
$ws.Range("H11").Value = 0.8344480000000001
$ws.Range("I11").Value = 0.001040856580779521
$ws.Range("J11").Value = 0.001040856580779521
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3487728915577651
$ws.Range("P11").Value = 0.3487728915577651
$ws.Range("Q11").Value = 46.95152476359824
$ws.Range("R11").Value = 422.5637228723841
$ws.Range("S11").Value = 0.0003630225593754021
$ws.Range("T11").Value = 0.0003630225593754021
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.2781493333333334
$ws.Range("H12").Value = 0.8344480000000001
$ws.Range("I12").Value = 0.001040856580779521
$ws.Range("J12").Value = 0.001040856580779521
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 68.09032333333333
$ws.Range("N12").Value = 204.27097
$ws.Range("O12").Value = 0.1406878008722904
$ws.Range("P12").Value = 0.1406878008722904
$ws.Range("Q12").Value = 18.93927804161778
$ws.Range("R12").Value = 170.45350237456
$ws.Range("S12").Value = 0.0001464358233733223
$ws.Range("T12").Value = 0.0001464358233733223
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.2781493333333334
$ws.Range("H13").Value = 0.8344480000000001
$ws.Range("I13").Value = 0.001040856580779521
$ws.Range("J13").Value = 0.001040856580779521
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 93.562673
$ws.Range("N13").Value = 280.688019
$ws.Range("O13").Value = 0.1933186106880956
$ws.Range("P13").Value = 0.1933186106880956
$ws.Range("Q13").Value = 26.02439511983467
$ws.Range("R13").Value = 234.219556078512
$ws.Range("S13").Value = 0.0002012169481218586
$ws.Range("T13").Value = 0.0002012169481218586
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 16.037618
$ws.Range("H14").Value = 48.112854
$ws.Range("I14").Value = 0.06001402209123193
$ws.Range("J14").Value = 0.06001402209123194
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 153.5290173333333
$ws.Range("N14").Value = 460.587052
$ws.Range("O14").Value = 0.3172206968818489
$ws.Range("P14").Value = 0.317220696881849
$ws.Range("Q14").Value = 2462.239731907378
$ws.Range("R14").Value = 22160.15758716641
$ws.Range("S14").Value = 0.01903768991046327
$ws.Range("T14").Value = 0.01903768991046327
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 16.037618
$ws.Range("H15").Value = 48.112854
$ws.Range("I15").Value = 0.06001402209123193
$ws.Range("J15").Value = 0.06001402209123194
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3487728915577651
$ws.Range("P15").Value = 0.3487728915577651
$ws.Range("Q15").Value = 2707.145149881581
$ws.Range("R15").Value = 24364.30634893423
$ws.Range("S15").Value = 0.02093126401877055
$ws.Range("T15").Value = 0.02093126401877056
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 16.037618
$ws.Range("H16").Value = 48.112854
$ws.Range("I16").Value = 0.06001402209123193
$ws.Range("J16").Value = 0.06001402209123194
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 68.09032333333333
$ws.Range("N16").Value = 204.27097
$ws.Range("O16").Value = 0.1406878008722904
$ws.Range("P16").Value = 0.1406878008722904
$ws.Range("Q16").Value = 1092.006595116487
$ws.Range("R16").Value = 9828.059356048379
$ws.Range("S16").Value = 0.008443240789516474
$ws.Range("T16").Value = 0.008443240789516476
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 16.037618
$ws.Range("H17").Value = 48.112854
$ws.Range("I17").Value = 0.06001402209123193
$ws.Range("J17").Value = 0.06001402209123194
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 93.562673
$ws.Range("N17").Value = 280.688019
$ws.Range("O17").Value = 0.1933186106880956
$ws.Range("P17").Value = 0.1933186106880956
$ws.Range("Q17").Value = 1500.522408632914
$ws.Range("R17").Value = 13504.70167769622
$ws.Range("S17").Value = 0.01160182737248164
$ws.Range("T17").Value = 0.01160182737248164
